$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): B11 4 -> 5, C11 -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): B12 60 -> 75, C12 -13 -> -15.6, E12 "47/112" -> "59.4/140"
$ws.Range("B12").Value = 75
$ws.Range("C12").Value = -15.6
$ws.Range("E12").Value = "59.4/140"
